# "Cambiado configuracion para mayor conf" - update the confidence scores
# (column I) on Sheet1 for every data row, then leave the selection where
# the author left it (K38). I41 holds =AVERAGEIF(I2:I39,"<>-100") and will
# auto-recalculate from the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$confidence = @{
    2  = 89.933999999999997
    3  = 85.524000000000001
    4  = 83.477999999999994
    5  = 89.578999999999994
    6  = 84.122
    7  = 89.242999999999995
    8  = 89.933999999999997
    9  = 85.524000000000001
    10 = 83.477999999999994
    11 = 89.578999999999994
    12 = 84.122
    13 = 89.242999999999995
    14 = 89.933999999999997
    15 = 85.524000000000001
    16 = 83.477999999999994
    17 = 89.578999999999994
    18 = 84.122
    19 = 89.242999999999995
    20 = 89.933999999999997
    21 = 85.524000000000001
    22 = 83.477999999999994
    23 = 89.578999999999994
    24 = 83.477999999999994
    25 = 89.578999999999994
    26 = 84.122
    27 = 89.242999999999995
    28 = 89.933999999999997
    29 = 85.524000000000001
    30 = 83.477999999999994
    31 = 89.578999999999994
    32 = 83.477999999999994
    33 = 89.578999999999994
    34 = 84.122
    35 = 89.242999999999995
    36 = 89.933999999999997
    37 = 85.524000000000001
    38 = 83.477999999999994
    39 = 89.578999999999994
}

foreach ($row in $confidence.Keys) {
    $ws.Cells.Item($row, 9).Value = $confidence[$row]
}

# Match the author's final on-screen selection for this sheet.
$ws.Range("K38").Select()
